# Generate Report for Handoff
# This script updates the "f9878666-90b1-48cc-b4fa-4a6a87ee1180.md" file's row
# (row 3) across the Overview, zh-cn and de-de sheets to reflect that the file
# is now "Ready for handoff" instead of "Handed back: in sync with en-US",
# along with updated handoff timestamps and (for the localized sheets) an
# error detail message about the handback file being stale, plus widening
# the "Error Detail" column so the message is readable.

$wb = $excel.ActiveWorkbook

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Ready for handoff"
$wsOverview.Range("F3").Value = "Ready for handoff"
$wsOverview.Range("G3").Value = "2016-08-25 22:47:31"

# ---- zh-cn sheet ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = "Ready for handoff"
$wsZhCn.Range("H3").Value = "2016-08-25 22:47:27"
$wsZhCn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b93447977b2c48c36e981273da778a442f6ed658/e2e/f9878666-90b1-48cc-b4fa-4a6a87ee1180.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1462d31bc92c5f164e2add9abc121bba10d13235/e2e/f9878666-90b1-48cc-b4fa-4a6a87ee1180.md."
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# ---- de-de sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = "Ready for handoff"
$wsDeDe.Range("H3").Value = "2016-08-25 22:47:31"
$wsDeDe.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/b93447977b2c48c36e981273da778a442f6ed658/e2e/f9878666-90b1-48cc-b4fa-4a6a87ee1180.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/1462d31bc92c5f164e2add9abc121bba10d13235/e2e/f9878666-90b1-48cc-b4fa-4a6a87ee1180.md."
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
